# Insert a new price record as row 82 (pushing the existing rows 82:151
# down to 83:152), matching the weekly data refresh described in the
# commit message ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 82:151 down to 83:152, inheriting formatting (incl. the
# date style on column D) from the row above, same as Excel's native
# "Insert Sheet Rows" behavior.
$ws.Rows(82).Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(82, 1).Value  = 7
$ws.Cells.Item(82, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(82, 3).Value  = "Ñuble"
$ws.Cells.Item(82, 4).Value  = 44669
$ws.Cells.Item(82, 5).Value  = 16
$ws.Cells.Item(82, 6).Value  = 100112045
$ws.Cells.Item(82, 7).Value  = "Zapallo"
$ws.Cells.Item(82, 8).Value  = "Camote"
$ws.Cells.Item(82, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(82, 10).Value = 200
$ws.Cells.Item(82, 11).Value = 300
$ws.Cells.Item(82, 12).Value = 350
$ws.Cells.Item(82, 13).Value = 325
$ws.Cells.Item(82, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(82, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(82, 16).Value = 325
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = "Hortaliza"
